{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the text replacements described by the commit diff:\n//   - document number \"11\" -> \"14\"\n//   - customer info block filled in (Megrendel\u0151 / C\u00edm / El\u00e9rhet\u0151s\u00e9g)\n//   - \"Megjegyz\u00e9s:\" comment appended\n//   - device info block filled in (Megnevez\u00e9s / T\u00edpus / Modell)\n//   - fault block filled in (Hibajelens\u00e9g / Tartoz\u00e9kok)\n//   - diagnosis filled in (Szerviz diagn\u00f3zis)\n//   - dates \"2024.07.24\" -> \"2024.07.28\" (all occurrences)\n\nconst body = context.document.body;\n\n// Simple \"find exactly one match and replace its text\" helper. Using\n// matchCase (and matchWholeWord where useful) keeps each search scoped to\n// a single, unambiguous hit so formatting on the run is preserved.\nasync function replaceOnce(findText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(findText, opts);\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${findText}\", found ${results.items.length}`\n    );\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// Replace every match (used for the date, which appears three times).\nasync function replaceAll(findText, newText, options) {\n  const opts = Object.assign({ matchCase: true }, options || {});\n  const results = body.search(findText, opts);\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n// 1) Document number, top right of the receipt header (\"\u00c1TV\u00c9ELI\n//    ELISMERV\u00c9NY ... 11\" -> \"... 14\"). Whole-word match avoids the \"11\"\n//    substring inside \"Rumi \u00fat 311.\".\nawait replaceOnce(\"11\", \"14\", { matchWholeWord: true });\n\n// 2) Customer block.\nawait replaceOnce(\"Megrendel\u0151: kuuuu\", \"Megrendel\u0151: P\u00e9lda P\u00e9ter\");\nawait replaceOnce(\"C\u00edm: sdsdsdsd\", \"C\u00edm: 9700 Szombathely Neml\u00e9tezik utca. 3\");\nawait replaceOnce(\n  \"El\u00e9rhet\u0151s\u00e9g: telefon  \",\n  \"El\u00e9rhet\u0151s\u00e9g: telefon  06301234567\"\n);\n\n// 3) Comment field.\nawait replaceOnce(\"Megjegyz\u00e9s:\", \"Megjegyz\u00e9s:uuheuhehu\");\n\n// 4) Device block.\nawait replaceOnce(\"Megnevez\u00e9s: \", \"Megnevez\u00e9s: F\u00fcnyiro\");\nawait replaceOnce(\"T\u00edpus: \", \"T\u00edpus: Ferrari\");\nawait replaceOnce(\"Modell: \", \"Modell: Igen\");\n\n// 5) Fault block.\nawait replaceOnce(\"Hibajelens\u00e9g: \", \"Hibajelens\u00e9g: Van\");\nawait replaceOnce(\"Tartoz\u00e9kok: \", \"Tartoz\u00e9kok: Nincs\");\n\n// 6) Diagnosis.\nawait replaceOnce(\"Szerviz diagn\u00f3zis: \", \"Szerviz diagn\u00f3zis: R\u00f3sz\");\n\n// 7) Dates (three occurrences across the document).\nawait replaceAll(\"2024.07.24\", \"2024.07.28\");\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the text replacements described by the commit diff:\n#   - document number \"11\" -> \"14\"\n#   - customer info block filled in (Megrendel\u0151 / C\u00edm / El\u00e9rhet\u0151s\u00e9g)\n#   - \"Megjegyz\u00e9s:\" comment appended\n#   - device info block filled in (Megnevez\u00e9s / T\u00edpus / Modell)\n#   - fault block filled in (Hibajelens\u00e9g / Tartoz\u00e9kok)\n#   - diagnosis filled in (Szerviz diagn\u00f3zis)\n#   - dates \"2024.07.24\" -> \"2024.07.28\" (all occurrences)\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2, wdFindContinue = 1 (search whole story, do not prompt)\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n# NOTE: named parameter binding (-FindText \"...\") is not reliable in this\n# interpreter, so the helper takes plain positional arguments instead.\nfunction Replace-Text($FindText, $ReplaceText, $MatchWholeWord = $false) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.Execute($FindText, $true, $MatchWholeWord, $false, $false, $false, $true, $wdFindContinue, $false, $ReplaceText, $wdReplaceAll)\n}\n\n# 1) Document number, top right of the receipt header (\"\u00c1TV\u00c9ELI\n#    ELISMERV\u00c9NY ... 11\" -> \"... 14\"). Whole-word match avoids the \"11\"\n#    substring inside \"Rumi \u00fat 311.\".\nReplace-Text \"11\" \"14\" $true\n\n# 2) Customer block.\nReplace-Text \"Megrendel\u0151: kuuuu\" \"Megrendel\u0151: P\u00e9lda P\u00e9ter\"\nReplace-Text \"C\u00edm: sdsdsdsd\" \"C\u00edm: 9700 Szombathely Neml\u00e9tezik utca. 3\"\nReplace-Text \"El\u00e9rhet\u0151s\u00e9g: telefon  \" \"El\u00e9rhet\u0151s\u00e9g: telefon  06301234567\"\n\n# 3) Comment field.\nReplace-Text \"Megjegyz\u00e9s:\" \"Megjegyz\u00e9s:uuheuhehu\"\n\n# 4) Device block.\nReplace-Text \"Megnevez\u00e9s: \" \"Megnevez\u00e9s: F\u00fcnyiro\"\nReplace-Text \"T\u00edpus: \" \"T\u00edpus: Ferrari\"\nReplace-Text \"Modell: \" \"Modell: Igen\"\n\n# 5) Fault block.\nReplace-Text \"Hibajelens\u00e9g: \" \"Hibajelens\u00e9g: Van\"\nReplace-Text \"Tartoz\u00e9kok: \" \"Tartoz\u00e9kok: Nincs\"\n\n# 6) Diagnosis.\nReplace-Text \"Szerviz diagn\u00f3zis: \" \"Szerviz diagn\u00f3zis: R\u00f3sz\"\n\n# 7) Dates (three occurrences across the document, wdReplaceAll covers all of them).\nReplace-Text \"2024.07.24\" \"2024.07.28\"\n"}
